# Relabel the K-Means "Cluster" column (column E) values according to the
# mapping produced by the bugfix: 0->4, 1->0, 2->1, 3->3 (unchanged), 4->2
# Column F (Cluster Centroid label) is intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{ 0 = 4; 1 = 0; 2 = 1; 3 = 3; 4 = 2 }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = [int]$cell.Value()
    $cell.Value = $map[$old]
}
